$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"/"In Advance" block gets
# pushed right by one), mirroring Excel's "Insert Sheet Columns" which takes
# its formatting from the column immediately to the left (M).
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab with S7 selected.
$ws.Activate() | Out-Null
$ws.Range("S7").Select() | Out-Null
